$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows to append (Serie date string, BCP, BCU, Otros)
$newRows = @(
    @("03-09-2021", 906, 5156, 77),
    @("04-09-2021", 906, 5156, 77),
    @("05-09-2021", 906, 5156, 77),
    @("06-09-2021", 908, 5168, 77)
)

$startRow = 247
$scratchRow = 300

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a date-like label (e.g. "03-09-2021") that Excel's
    # Value setter would otherwise auto-convert into a real date serial.
    # Build it as a text formula result in a scratch cell, then copy/paste
    # it as a value so it lands as plain text (matching existing cells)
    # without picking up a new number-format style.
    $scratchCell = $ws.Cells.Item($scratchRow, 1)
    $scratchCell.Formula = "=""" + $row[0] + """"
    $scratchCell.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Cells.Item($scratchRow, 1).Clear()
$excel.CutCopyMode = $false
